$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

$ws.Range("A9").Value = "Viki"
$ws.Range("B2").Copy()
$ws.Range("B9").PasteSpecial(-4122)
$ws.Range("B9").Value = 45321
$ws.Range("C9").Formula = "=8+45/60"
$ws.Range("D9").Formula = "=9+35/60"

$ws.Range("A10").Value = "Aris"
$ws.Range("B2").Copy()
$ws.Range("B10").PasteSpecial(-4122)
$ws.Range("B10").Value = 45322
$ws.Range("C10").Formula = "=8+45/60"
$ws.Range("D10").Formula = "=9+35/60"
$ws.Range("F10").Value = "Items"

$ws.Range("F9").Value = "Storyline"

$null = $ws.Range("F10").Select()
